$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row at row 8 (pushes existing rows 8:37 down to 9:38)
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the "RUN SCRIPT" web service URI entry.
# Column A: command name -- reuse the existing "-- not implemented --" text
# (leading apostrophe forces the quote-prefix text style, like other "--" rows).
$ws.Range("A8").Value = "'-- not implemented --"

# Column B: method label, computed the same way as the rest of the column.
$ws.Range("B8").Formula = "=LEFT(A8,SEARCH("" "",A8)-1)"

# Column C: HTTP verb.
$ws.Range("C8").Value = "POST"

# Column D: URI template for running a script directly without a job.
$ws.Range("D8").Value = "scripts/-/running/{scriptName}"

# Update the hidden _FilterDatabase defined name so it spans the new last row.
# (Iterate by index -- looking it up by bare name after the row insert
# returns a stale/blank RefersTo in this runtime.)
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $definedName = $wb.Names.Item($i)
    if ($definedName.Name -like "*_FilterDatabase*") {
        $definedName.RefersTo = "=Commands!`$B`$1:`$B`$38"
    }
}

# Restore the active selection to A2 (top-left of the frozen scrollable area).
$ws.Range("A2").Select() | Out-Null
